$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that contain purely numeric-looking text must be forced
# to remain text (matching the original inlineStr type) otherwise Excel
# auto-converts them to numbers and formatting (e.g. trailing zeros) is lost.

$ws.Range("D2").Value = "51.462.74"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "2.930.15"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "377.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.541"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.587"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0838"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("D14").Value = "3.396.07"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.65%  "
$ws.Range("D16").Value = "2.932.11"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.944"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.91%  "
$ws.Range("D18").Value = "51.476.45"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").Value = "0.0₃0950"
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +17.62%  "
$ws.Range("E27").Value = "  -4.63%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -4.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.102"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "52.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "34.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.14%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0428"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E41").Value = "  -6.60%  "
$ws.Range("E42").Value = "  -5.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.114"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.72%  "
$ws.Range("E46").Value = "  -4.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.275"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +14.92%  "
$ws.Range("D48").Value = "2.019.46"
$ws.Range("E48").Value = "  -4.34%  "
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0321"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.35%  "
